$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "car/car070.png"
$ws.Range("C2").Value = "fesseln"
$ws.Range("D2").Value = "car"
$ws.Range("B3").Value = "dog/dog112.png"
$ws.Range("C3").Value = "klappen"
$ws.Range("D3").Value = "dog"
$ws.Range("B4").Value = "car/car079.png"
$ws.Range("C4").Value = "liefern"
$ws.Range("D4").Value = "car"
$ws.Range("B5").Value = "car/car107.png"
$ws.Range("C5").Value = "füllen"
$ws.Range("D5").Value = "car"
$ws.Range("B6").Value = "car/car118.png"
$ws.Range("C6").Value = "stechen"
$ws.Range("D6").Value = "car"
$ws.Range("B7").Value = "dog/dog119.png"
$ws.Range("C7").Value = "formen"
$ws.Range("D7").Value = "dog"
$ws.Range("B8").Value = "dog/dog093.png"
$ws.Range("C8").Value = "opfern"
$ws.Range("D8").Value = "dog"
$ws.Range("B9").Value = "dog/dog113.png"
$ws.Range("C9").Value = "fliegen"
$ws.Range("D9").Value = "dog"
$ws.Range("B10").Value = "car/car093.png"
$ws.Range("C10").Value = "enden"
$ws.Range("D10").Value = "car"
$ws.Range("B11").Value = "dog/dog103.png"
$ws.Range("C11").Value = "jubeln"
$ws.Range("D11").Value = "dog"
$ws.Range("B12").Value = "dog/dog067.png"
$ws.Range("C12").Value = "krachen"
$ws.Range("D12").Value = "dog"
$ws.Range("B13").Value = "car/car099.png"
$ws.Range("C13").Value = "nehmen"
$ws.Range("D13").Value = "car"
$ws.Range("B14").Value = "car/car064.png"
$ws.Range("C14").Value = "regnen"
$ws.Range("D14").Value = "car"
$ws.Range("B15").Value = "car/car108.png"
$ws.Range("C15").Value = "sieben"
$ws.Range("D15").Value = "car"
$ws.Range("B16").Value = "dog/dog091.png"
$ws.Range("C16").Value = "posten"
$ws.Range("D16").Value = "dog"
$ws.Range("B17").Value = "car/car115.png"
$ws.Range("C17").Value = "dauern"
$ws.Range("D17").Value = "car"
$ws.Range("B18").Value = "car/car066.png"
$ws.Range("C18").Value = "scheitern"
$ws.Range("D18").Value = "car"
$ws.Range("B19").Value = "dog/dog108.png"
$ws.Range("C19").Value = "bitten"
$ws.Range("D19").Value = "dog"
$ws.Range("B20").Value = "dog/dog089.png"
$ws.Range("C20").Value = "laufen"
$ws.Range("D20").Value = "dog"
$ws.Range("B21").Value = "dog/dog098.png"
$ws.Range("C21").Value = "biegen"
$ws.Range("D21").Value = "dog"
$ws.Range("B22").Value = "dog/dog102.png"
$ws.Range("C22").Value = "tauschen"
$ws.Range("D22").Value = "dog"
$ws.Range("B23").Value = "dog/dog115.png"
$ws.Range("C23").Value = "backen"
$ws.Range("D23").Value = "dog"
$ws.Range("B24").Value = "car/car088.png"
$ws.Range("C24").Value = "wenden"
$ws.Range("D24").Value = "car"
$ws.Range("B25").Value = "car/car105.png"
$ws.Range("C25").Value = "schmecken"
$ws.Range("D25").Value = "car"
$ws.Range("B26").Value = "car/car078.png"
$ws.Range("C26").Value = "kaufen"
$ws.Range("D26").Value = "car"
$ws.Range("B27").Value = "car/car083.png"
$ws.Range("C27").Value = "saufen"
$ws.Range("D27").Value = "car"
$ws.Range("B28").Value = "dog/dog082.png"
$ws.Range("C28").Value = "gelten"
$ws.Range("D28").Value = "dog"
$ws.Range("B29").Value = "dog/dog094.png"
$ws.Range("C29").Value = "hupen"
$ws.Range("D29").Value = "dog"
$ws.Range("B30").Value = "dog/dog078.png"
$ws.Range("C30").Value = "bleiben"
$ws.Range("D30").Value = "dog"
$ws.Range("B31").Value = "dog/dog081.png"
$ws.Range("C31").Value = "pflegen"
$ws.Range("D31").Value = "dog"
$ws.Range("B32").Value = "car/car119.png"
$ws.Range("C32").Value = "langen"
$ws.Range("D32").Value = "car"
$ws.Range("B33").Value = "car/car085.png"
$ws.Range("C33").Value = "fliehen"
$ws.Range("D33").Value = "car"
